$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Recorded By" (column G) value needs its first two
# comma-separated entries swapped (a trailing 3rd entry, if any, stays put).
$rowsToSwap = @(2,4,5,7,8,11,17,28,30,31,33,34,37,43,54,56,57,59,60,63,69,80,81,82,87,93,94,96,106,107,108,113,119,120,122,132,133,134,139,145,146,148)

foreach ($r in $rowsToSwap) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    $parts = $val -split ', '
    if ($parts.Length -ge 2) {
        $tmp = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $tmp
    }
    $cell.Value = [string]::Join(', ', $parts)
}
